$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows for A2:C16 (player, position, team) after the roster update.
$data = @(
    @("Donovan Mitchell", "PG,SG", "Cleveland Cavaliers"),
    @("Jamal Murray", "PG,SG", "Denver Nuggets"),
    @("Malik Beasley", "SG", "Detroit Pistons"),
    @("Dyson Daniels", "PG,SG", "Atlanta Hawks"),
    @("Jaden McDaniels", "SF,PF", "Minnesota Timberwolves"),
    @("Naz Reid", "PF,C", "Minnesota Timberwolves"),
    @("Alexandre Sarr", "PF,C", "Washington Wizards"),
    @("Toumani Camara", "SF,PF", "Portland Trail Blazers"),
    @("Myles Turner", "C", "Indiana Pacers"),
    @("Jonas Valanciunas", "C", "Washington Wizards"),
    @("Domantas Sabonis", "C", "Sacramento Kings"),
    @("Victor Wembanyama", "C", "San Antonio Spurs"),
    @("Josh Hart", "SF,PF", "New York Knicks"),
    @("Tari Eason", "SF,PF", "Houston Rockets"),
    @("De'Andre Hunter", "SF,PF", "Atlanta Hawks")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row++
}
